$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated hand-count values for B244:B264
$values = @(
    95,
    87,
    91,
    80,
    87,
    58,
    55,
    29,
    33,
    50,
    67,
    24,
    51,
    65,
    31,
    28,
    36,
    57,
    18,
    48,
    16
)

$startRow = 244
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the active selection on the sheet to match the new state
$ws.Activate()
$ws.Range("B2:B264").Select()
